$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$wCell = $ws.Range("W5")
$wCell.Formula = '="2506"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 6
$wCell = $ws.Range("W6")
$wCell.Formula = '="2720"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 18
$ws.Range("V18").Value = 19
$ws.Range("V18").Interior.Color = 65535
$wCell = $ws.Range("W18")
$wCell.Formula = '="3995"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 19
$wCell = $ws.Range("W19")
$wCell.Formula = '="2636"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 20
$ws.Range("V20").Value = 28
$wCell = $ws.Range("W20")
$wCell.Formula = '="4213"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 21
$ws.Range("V21").Value = 37
$wCell = $ws.Range("W21")
$wCell.Formula = '="4708"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 22
$ws.Range("V22").Value = 20
$wCell = $ws.Range("W22")
$wCell.Formula = '="4681"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 23
$wCell = $ws.Range("W23")
$wCell.Formula = '="5283"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 24
$wCell = $ws.Range("W24")
$wCell.Formula = '="4832"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 25
$wCell = $ws.Range("W25")
$wCell.Formula = '="5118"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 26
$ws.Range("V26").Value = 0
$ws.Range("V26").Interior.Color = 255

# Row 29
$ws.Range("V29").Value = 16
$ws.Range("V29").Interior.Color = 65535
$wCell = $ws.Range("W29")
$wCell.Formula = '="3050"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 30
$ws.Range("V30").Value = 33
$ws.Range("V30").Interior.Color = 32768
$wCell = $ws.Range("W30")
$wCell.Formula = '="4956"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 31
$ws.Range("V31").Value = 33
$wCell = $ws.Range("W31")
$wCell.Formula = '="4862"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 32
$ws.Range("V32").Value = 2
$ws.Range("V32").Interior.Color = 65535
$wCell = $ws.Range("W32")
$wCell.Formula = '="2696"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 35
$wCell = $ws.Range("W35")
$wCell.Formula = '="4601"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 36
$wCell = $ws.Range("W36")
$wCell.Formula = '="2727"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 37
$wCell = $ws.Range("W37")
$wCell.Formula = '="4541"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 38
$ws.Range("V38").Value = 0
$ws.Range("V38").Interior.Color = 255
$wCell = $ws.Range("W38")
$wCell.Formula = '="4918"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 39
$wCell = $ws.Range("W39")
$wCell.Formula = '="4575"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 41
$wCell = $ws.Range("W41")
$wCell.Formula = '="4281"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 42
$wCell = $ws.Range("W42")
$wCell.Formula = '="3044"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 44
$ws.Range("W44").Value = 4860

# Row 45
$ws.Range("V45").Value = 29
$wCell = $ws.Range("W45")
$wCell.Formula = '="3987"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 46
$wCell = $ws.Range("W46")
$wCell.Formula = '="3915"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 47
$ws.Range("V47").Value = 33
$wCell = $ws.Range("W47")
$wCell.Formula = '="5226"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 48
$ws.Range("V48").Value = 30
$ws.Range("V48").Interior.Color = 16777215
$wCell = $ws.Range("W48")
$wCell.Formula = '="4817"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 49
$wCell = $ws.Range("W49")
$wCell.Formula = '="4695"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 50
$wCell = $ws.Range("W50")
$wCell.Formula = '="4777"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 51
$ws.Range("V51").Value = 16
$wCell = $ws.Range("W51")
$wCell.Formula = '="3875"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 52
$wCell = $ws.Range("W52")
$wCell.Formula = '="4949"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 53
$ws.Range("V53").Value = 8
$wCell = $ws.Range("W53")
$wCell.Formula = '="3635"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 54
$wCell = $ws.Range("W54")
$wCell.Formula = '="4623"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 55
$ws.Range("V55").Value = 10
$wCell = $ws.Range("W55")
$wCell.Formula = '="3705"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 56
$ws.Range("V56").Value = 30
$wCell = $ws.Range("W56")
$wCell.Formula = '="5165"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 57
$ws.Range("V57").Value = 18
$ws.Range("V57").Interior.Color = 65535
$wCell = $ws.Range("W57")
$wCell.Formula = '="4203"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 58
$wCell = $ws.Range("W58")
$wCell.Formula = '="4106"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 59
$ws.Range("V59").Value = 18
$ws.Range("V59").Interior.Color = 65535
$wCell = $ws.Range("W59")
$wCell.Formula = '="4074"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 60
$wCell = $ws.Range("W60")
$wCell.Formula = '="4206"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 62
$wCell = $ws.Range("W62")
$wCell.Formula = '="3988"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 63
$ws.Range("V63").Value = 22
$wCell = $ws.Range("W63")
$wCell.Formula = '="4035"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 64
$ws.Range("V64").Value = 11
$ws.Range("V64").Interior.Color = 65535
$wCell = $ws.Range("W64")
$wCell.Formula = '="4143"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 65
$ws.Range("V65").Value = 23
$wCell = $ws.Range("W65")
$wCell.Formula = '="3888"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 69
$wCell = $ws.Range("W69")
$wCell.Formula = '="2956"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 70
$wCell = $ws.Range("W70")
$wCell.Formula = '="1518"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 73
$wCell = $ws.Range("W73")
$wCell.Formula = '="2605"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 76
$ws.Range("W76").Value = 4114

# Row 77
$wCell = $ws.Range("W77")
$wCell.Formula = '="3675"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 78
$wCell = $ws.Range("W78")
$wCell.Formula = '="2862"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 79
$wCell = $ws.Range("W79")
$wCell.Formula = '="1306"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 81
$wCell = $ws.Range("W81")
$wCell.Formula = '="2647"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 84
$wCell = $ws.Range("W84")
$wCell.Formula = '="1524"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 96
$wCell = $ws.Range("W96")
$wCell.Formula = '="2790"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 115
$ws.Range("W115").Value = 6011

# Row 116
$ws.Range("W116").Value = 4944

# Row 117
$ws.Range("W117").Value = 4672

# Row 118
$ws.Range("V118").Value = 20
$wCell = $ws.Range("W118")
$wCell.Formula = '="5403"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)

# Row 119
$ws.Range("V119").Value = 0
$ws.Range("V119").Interior.Color = 255
$wCell = $ws.Range("W119")
$wCell.Formula = '="1609"'
$wCell.Copy()
$wCell.PasteSpecial(-4163)
